# 发货模板.xlsx — update instructions text in A1 and extend the merged
# header row's styled cells, update selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Rebuild the rich-text instructions in A1.
#    Run 1: "注：导入时请删除该行" (now left with default/no explicit run
#    formatting -- it simply inherits the cell's own style).
#    Run 2: the numbered usage notes, with a new item 2 inserted about
#    order-number cell format, renumbering the remaining items.
# ---------------------------------------------------------------------
$note = "注：导入时请删除该行"
$rules = @"

   1.必填项：订单号、物流公司、快递单号、快递时间；
   2.关于订单号单元格格式，必须设置为文本格式
   3.商品名称：非必填，如果不填写，则默认订单下所有商品全部发货；如果填写商品名称，则只发该商品对应的发货数量；
   4.发货数量：非必填，如果不填写，则默认全部发货；如果填写，则只发填写的数量；发货数量填写后商品名称必填；
   5.快递时间：必须为“年-月-日 时:分:秒”格式；比如：2022-10-10 10:00:00；不能写成2022/10/10 10:00:00或20221010 10:00:00；

"@

$a1 = $ws.Range("A1")
$a1.Value = $note + $rules

# Re-apply the explanatory-notes formatting to the second run only, so the
# leading "注：..." fragment is left with no run-level formatting overrides.
$run2 = $a1.Characters($note.Length + 1, $rules.Length)
$run2.Font.Size = 11
$run2.Font.Name = "宋体"
$run2.Font.Color = 0

# ---------------------------------------------------------------------
# 2. The merged A1:F1 header cell now carries explicit (empty) styled
#    cells across the whole merge range -- copy A1's formatting onto
#    B1:F1.
# ---------------------------------------------------------------------
$a1.Copy() | Out-Null
$ws.Range("B1:F1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3. Move the active selection to A2.
# ---------------------------------------------------------------------
$ws.Range("A2").Select() | Out-Null
